# "Trying to handle secure and non-secure site"
# Upgrade the CHROME / MC.Browser "Global" rows to the https:// URL while
# leaving the FIREFOX row on the plain (non-secure) host, and add the new
# secure URL string to the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Global")

$secureUrl = "https://advantageonlineshopping.com/"

# Row 2: Browser / CHROME -> secure URL
$ws.Range("C2").Value = $secureUrl

# Row 6: Device / MC.Browser / Android -> secure URL
$ws.Range("C6").Value = $secureUrl

# Row 3 (Browser / FIREFOX) is intentionally left on the non-secure URL.

# Widen column C (best-fit style column) so the longer URL text keeps fitting.
$ws.Columns.Item(3).ColumnWidth = 31.928385416666664

# Leave the active selection on C3, matching where the author ended up.
$ws.Range("C3").Select()

# "Action 1" remains the active tab, as in the original file.
$wb.Worksheets.Item("Action 1").Activate()
